# Time estimation workbook update:
#  - add a "Buffer:" row (row 28) below the Tabelle3 table on sheet "UC"
#  - buffer factor of 2.5, buffer time = total time * buffer FP factor
#  - expand the Tabelle3 table/autofilter range to include the new row
#  - bold the new computed buffer-time cell
#  - move the active selection/view down to the newly added row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC")
$ws.Activate()

# New data row right below the existing "Tabelle3" table (A23:C27 -> A23:C28)
$ws.Range("A28").Value = "Buffer:"
$ws.Range("B28").Value = 2.5
$ws.Range("C28").Formula = "=C27*Tabelle3[[#This Row],[FP]]"

# Make the total buffer time stand out in bold, like the other summary rows
$ws.Range("C28").Font.Bold = $true

# Grow the table (and its autofilter) so it covers the freshly added row
$lo = $ws.ListObjects.Item("Tabelle3")
$lo.Resize($ws.Range("A23:C28"))

# Scroll / select near the new row, matching where the editor ended up
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("B29").Select()
